$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Add a new bold paragraph ("Play Fenix Play Deluxe for Free -
#    Exciting Firebird Slot Game") right before the last paragraph
#    (the italic "feature image" prompt paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

# Copy a plain, unformatted paragraph so the freshly split paragraph
# mark doesn't inherit the italic formatting that is sitting right
# next to the insertion point.
$template = $d.Paragraphs(3)
$template.Range.Copy()

$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.Paste()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newCount - 1)
$newPara.Range.Text = "Play Fenix Play Deluxe for Free - Exciting Firebird Slot Game"

$newPara = $d.Paragraphs($newCount - 1)
$newParaTextOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaTextOnly.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Swap the "Create a cartoon-style feature image..." image-prompt
#    text for the meta-description copy, keeping the italic run.
# ------------------------------------------------------------------
$oldText = "Create a cartoon-style feature image for Fenix Play Deluxe that features a happy Maya warrior wearing glasses. The image should be colorful and eye-catching, drawing the attention of potential players. The Maya warrior should be holding a staff, standing in front of a backdrop that includes the ancient firebird rising from the ashes. The image should be fun and lighthearted, capturing the excitement and adventure of the game. The overall design should be simple and clean, with bold and bright colors that pop. The Maya warrior should be the central focus of the image, with the firebird in the background to highlight the theme of the game."
$newText = "Read our review of Fenix Play Deluxe - a classic three-reel, five-pay line game with engaging bonus features. Play for free today."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
